$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.026.59"

$ws.Range("E2").Value = "  -0.87%  "

$ws.Range("D3").Value = "2.247.62"

$ws.Range("E3").Value = "  -1.38%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").NumberFormat = "General"

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.05"
$ws.Range("D5").NumberFormat = "General"

$ws.Range("E5").Value = "  -0.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.48"
$ws.Range("D6").NumberFormat = "General"

$ws.Range("E6").Value = "  -1.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.523"
$ws.Range("D7").NumberFormat = "General"

$ws.Range("E7").Value = "  -1.40%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.487"
$ws.Range("D9").NumberFormat = "General"

$ws.Range("E9").Value = "  -1.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.69"
$ws.Range("D10").NumberFormat = "General"

$ws.Range("E10").Value = "  -3.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0805"
$ws.Range("D11").NumberFormat = "General"

$ws.Range("E11").Value = "  +1.02%  "

$ws.Range("E12").Value = "  +1.15%  "

$ws.Range("E13").Value = "  +0.91%  "

$ws.Range("D14").Value = "2.594.03"

$ws.Range("E14").Value = "  -1.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.40"
$ws.Range("D15").NumberFormat = "General"

$ws.Range("E15").Value = "  -0.61%  "

$ws.Range("D16").Value = "2.249.68"

$ws.Range("E16").Value = "  -1.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.778"
$ws.Range("D17").NumberFormat = "General"

$ws.Range("E17").Value = "  -3.08%  "

$ws.Range("D18").Value = "41.902.51"

$ws.Range("E18").Value = "  -0.98%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.12"
$ws.Range("D19").NumberFormat = "General"

$ws.Range("D20").Value = "0.0₃0900"

$ws.Range("E20").Value = "  -1.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.91"
$ws.Range("D21").NumberFormat = "General"

$ws.Range("E21").Value = "  -1.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.06"
$ws.Range("D22").NumberFormat = "General"

$ws.Range("E22").Value = "  -0.97%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.42"
$ws.Range("D23").NumberFormat = "General"

$ws.Range("E23").Value = "  -2.35%  "

$ws.Range("E24").Value = "  -1.40%  "

$ws.Range("E25").Value = "  -0.42%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.26"
$ws.Range("D27").NumberFormat = "General"

$ws.Range("E27").Value = "  -2.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.61"
$ws.Range("D28").NumberFormat = "General"

$ws.Range("E28").Value = "  -0.08%  "

$ws.Range("E29").Value = "  +0.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.46"
$ws.Range("D30").NumberFormat = "General"

$ws.Range("E30").Value = "  -0.64%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.27"
$ws.Range("D31").NumberFormat = "General"

$ws.Range("E31").Value = "  +4.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").NumberFormat = "General"

$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.14"
$ws.Range("D33").NumberFormat = "General"

$ws.Range("E33").Value = "  -2.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.08"
$ws.Range("D34").NumberFormat = "General"

$ws.Range("E34").Value = "  -2.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.42"
$ws.Range("D35").NumberFormat = "General"

$ws.Range("E35").Value = "  +2.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0717"
$ws.Range("D36").NumberFormat = "General"

$ws.Range("E36").Value = "  -3.28%  "

$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("E38").Value = "  -0.43%  "

$ws.Range("E39").Value = "  -3.37%  "

$ws.Range("E40").Value = "  -3.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.04"
$ws.Range("D41").NumberFormat = "General"

$ws.Range("E41").Value = "  -2.07%  "

$ws.Range("D42").Value = "1.940.16"

$ws.Range("E42").Value = "  -2.93%  "

$ws.Range("E43").Value = "  -2.14%  "

$ws.Range("B44").Value = "EnergySwap"

$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.53"
$ws.Range("D44").NumberFormat = "General"

$ws.Range("E44").Value = "  -2.45%  "

$ws.Range("B45").Value = "ApeXProtocol"

$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.18"
$ws.Range("D45").NumberFormat = "General"

$ws.Range("E45").Value = "  -9.83%  "

$ws.Range("E46").Value = "  -3.34%  "

$ws.Range("E47").Value = "  -3.75%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.42"
$ws.Range("D48").NumberFormat = "General"

$ws.Range("E48").Value = "  +0.85%  "

$ws.Range("D49").Value = "2.466.16"

$ws.Range("E49").Value = "  -1.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.18"
$ws.Range("D50").NumberFormat = "General"

$ws.Range("E50").Value = "  -1.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "91.02"
$ws.Range("D51").NumberFormat = "General"

$ws.Range("E51").Value = "  -1.15%  "
